# Add read and write scripts
# - sheet1 (data_types): remove the now-unused leading blank column/row so
#   the table starts at A1 instead of B2, and relabel "boolan" -> "boolean".
# - sheet2 (people): add "age" and "city" columns with sample data.
# - sheet3 (expenses): untouched (string-table renumbering only, handled
#   automatically by the engine).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 ("data_types"): shift the table from B2:D11 up-and-left to
# A1:C10 by deleting the now-empty column A and row 1.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(1).Delete()
$ws1.Rows.Item(1).Delete()

# ---------------------------------------------------------------------
# Sheet 2 ("people"): append "age" and "city" columns.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("D1").Value = "age"
$ws2.Range("E1").Value = "city"

$ages = @(12, 25, 4, 64, 45, 85, 55, 11)
$cities = @("London", "Paris", "Madrid", "New York", "Lisbon", "Shanghai", "Tokyo", "Camberra")

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 4).Value = $ages[$i]
    $ws2.Cells.Item($r, 5).Value = $cities[$i]
}

$ws2.Activate()
$ws2.Range("E10").Select()

# Fix the "boolan" typo -> "boolean" (done last so the new shared string
# lands at the end of the string table, after the age/city values)
$ws1.Range("A10").Value = "boolean"

$ws1.Activate()
$ws1.Range("C13").Select()
